# Auto-generated edit script: updates market-board derived price/profit
# columns (H-N) on several sheets, matching the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 6487
$ws.Range("J80").Value = 6499.5
$ws.Range("L80").Value = 19498.5
$ws.Range("N80").Value = -21494.5
$ws.Range("H83").Value = 6487
$ws.Range("J83").Value = 6499.5
$ws.Range("L83").Value = 58495.5
$ws.Range("N83").Value = -68479.5
$ws.Range("H86").Value = 4233.3335
$ws.Range("I86").Value = 3879.8
$ws.Range("K86").Value = 3879.8
$ws.Range("M86").Value = -2756.8
$ws.Range("H89").Value = 4233.3335
$ws.Range("I89").Value = 3879.8
$ws.Range("K89").Value = 19399
$ws.Range("M89").Value = -13783
$ws.Range("H121").Value = 3674.5
$ws.Range("J121").Value = 3674.5
$ws.Range("L121").Value = 11023.5
$ws.Range("N121").Value = -14517.5
$ws.Range("H132").Value = 3520.7585
$ws.Range("I132").Value = 3346.6667
$ws.Range("K132").Value = 10040.0001
$ws.Range("M132").Value = -7510.000100000001
$ws.Range("H137").Value = 418403.75
$ws.Range("I137").Value = 501392.1
$ws.Range("K137").Value = 1504176.3
$ws.Range("M137").Value = -1501626.3
$ws.Range("H138").Value = 3817.9119
$ws.Range("I138").Value = 3365.8845
$ws.Range("K138").Value = 10097.6535
$ws.Range("M138").Value = -4957.6535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3758.3
$ws.Range("I74").Value = 1499.1666
$ws.Range("K74").Value = 1499.1666
$ws.Range("M74").Value = -625.1666
$ws.Range("H77").Value = 3758.3
$ws.Range("I77").Value = 1499.1666
$ws.Range("K77").Value = 7495.833000000001
$ws.Range("M77").Value = -3127.833000000001
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H122").Value = 3434.7334
$ws.Range("I122").Value = 2690.125
$ws.Range("K122").Value = 8070.375
$ws.Range("M122").Value = -5620.375
$ws.Range("H125").Value = 101249.5
$ws.Range("J125").Value = 101249.5
$ws.Range("L125").Value = 101249.5
$ws.Range("N125").Value = -111089.5
$ws.Range("H132").Value = 1370.1538
$ws.Range("I132").Value = 1123.9143
$ws.Range("K132").Value = 3371.7429
$ws.Range("M132").Value = -841.7428999999997

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1918.6666
$ws.Range("I5").Value = 2102.6
$ws.Range("K5").Value = 2102.6
$ws.Range("M5").Value = -1989.6
$ws.Range("H95").Value = 20966.666
$ws.Range("J95").Value = 20966.666
$ws.Range("L95").Value = 20966.666
$ws.Range("N95").Value = -26458.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 63
$ws.Range("I7").Value = 67
$ws.Range("J7").Value = 51
$ws.Range("K7").Value = 67
$ws.Range("L7").Value = 51
$ws.Range("M7").Value = 46
$ws.Range("N7").Value = -277
$ws.Range("H31").Value = 2866.3809
$ws.Range("I31").Value = 1663.909
$ws.Range("K31").Value = 1663.909
$ws.Range("M31").Value = -1368.909
$ws.Range("H34").Value = 2866.3809
$ws.Range("I34").Value = 1663.909
$ws.Range("K34").Value = 1663.909
$ws.Range("M34").Value = -1461.909
$ws.Range("H58").Value = 5299.95
$ws.Range("I58").Value = 3210.0667
$ws.Range("K58").Value = 3210.0667
$ws.Range("M58").Value = -3007.0667
$ws.Range("H132").Value = 1817.8276
$ws.Range("I132").Value = 1757.75
$ws.Range("K132").Value = 5273.25
$ws.Range("M132").Value = -2743.25
$ws.Range("H134").Value = 4275.5483
$ws.Range("I134").Value = 3602.8635
$ws.Range("K134").Value = 10808.5905
$ws.Range("M134").Value = -8273.5905
$ws.Range("H136").Value = 5299.95
$ws.Range("I136").Value = 3210.0667
$ws.Range("K136").Value = 9630.2001
$ws.Range("M136").Value = -7080.2001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 389.83334
$ws.Range("I23").Value = 147
$ws.Range("J23").Value = 511.25
$ws.Range("K23").Value = 441
$ws.Range("L23").Value = 1533.75
$ws.Range("M23").Value = -206
$ws.Range("N23").Value = -2003.75
$ws.Range("H98").Value = 4248.926
$ws.Range("I98").Value = 4695.1665
$ws.Range("J98").Value = 4121.4287
$ws.Range("K98").Value = 14085.4995
$ws.Range("L98").Value = 12364.2861
$ws.Range("M98").Value = -12587.4995
$ws.Range("N98").Value = -15360.2861
$ws.Range("H103").Value = 2176.0908
$ws.Range("I103").Value = 361.33334
$ws.Range("J103").Value = 2856.625
$ws.Range("K103").Value = 1084.00002
$ws.Range("L103").Value = 8569.875
$ws.Range("M103").Value = -205.0000199999999
$ws.Range("N103").Value = -10327.875
$ws.Range("H132").Value = 2753.6
$ws.Range("I132").Value = 3054.5
$ws.Range("K132").Value = 27490.5
$ws.Range("M132").Value = -24960.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H70").Value = 5586.1924
$ws.Range("I70").Value = 5389.7
$ws.Range("K70").Value = 5389.7
$ws.Range("M70").Value = -5119.7
$ws.Range("H73").Value = 5586.1924
$ws.Range("I73").Value = 5389.7
$ws.Range("K73").Value = 5389.7
$ws.Range("M73").Value = -4453.7
$ws.Range("H92").Value = 9209
$ws.Range("J92").Value = 9209
$ws.Range("L92").Value = 9209
$ws.Range("N92").Value = -12953
$ws.Range("H107").Value = 296.75
$ws.Range("I107").Value = 297.7143
$ws.Range("K107").Value = 297.7143
$ws.Range("M107").Value = 1622.2857
$ws.Range("H132").Value = 3815.6553
$ws.Range("J132").Value = 2528.8333
$ws.Range("L132").Value = 7586.499899999999
$ws.Range("N132").Value = -12646.4999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1363.125
$ws.Range("I22").Value = 609.8
$ws.Range("J22").Value = 1705.5454
$ws.Range("K22").Value = 609.8
$ws.Range("L22").Value = 1705.5454
$ws.Range("M22").Value = -314.8
$ws.Range("N22").Value = -2295.5454
$ws.Range("H27").Value = 1363.125
$ws.Range("I27").Value = 609.8
$ws.Range("J27").Value = 1705.5454
$ws.Range("K27").Value = 609.8
$ws.Range("L27").Value = 1705.5454
$ws.Range("M27").Value = -502.8
$ws.Range("N27").Value = -1919.5454
$ws.Range("H55").Value = 5000469.5
$ws.Range("I55").Value = 5555926.5
$ws.Range("J55").Value = 1356.5
$ws.Range("K55").Value = 5555926.5
$ws.Range("L55").Value = 1356.5
$ws.Range("M55").Value = -5555753.5
$ws.Range("N55").Value = -1702.5
$ws.Range("H82").Value = 1914.1333
$ws.Range("I82").Value = 1515.8
$ws.Range("K82").Value = 1515.8
$ws.Range("M82").Value = -1154.8
$ws.Range("H85").Value = 1914.1333
$ws.Range("I85").Value = 1515.8
$ws.Range("K85").Value = 1515.8
$ws.Range("M85").Value = -267.8
$ws.Range("H100").Value = 1062
$ws.Range("I100").Value = 593
$ws.Range("K100").Value = 593
$ws.Range("M100").Value = -52
$ws.Range("H132").Value = 4799.357
$ws.Range("I132").Value = 3797.75
$ws.Range("K132").Value = 11393.25
$ws.Range("M132").Value = -8863.25
$ws.Range("H136").Value = 2048.8936
$ws.Range("I136").Value = 1317.7097
$ws.Range("J136").Value = 3465.5625
$ws.Range("K136").Value = 3953.1291
$ws.Range("L136").Value = 10396.6875
$ws.Range("M136").Value = -1403.1291
$ws.Range("N136").Value = -15496.6875
$ws.Range("H138").Value = 76857
$ws.Range("J138").Value = 76857
$ws.Range("L138").Value = 76857
$ws.Range("N138").Value = -87137
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 52999
$ws.Range("J70").Value = 52999
$ws.Range("L70").Value = 52999
$ws.Range("N70").Value = -53629
$ws.Range("H73").Value = 52999
$ws.Range("J73").Value = 52999
$ws.Range("L73").Value = 52999
$ws.Range("N73").Value = -55183
$ws.Range("H81").Value = 7392.1113
$ws.Range("I81").Value = 8224.286
$ws.Range("K81").Value = 16448.572
$ws.Range("M81").Value = -15387.572
$ws.Range("H84").Value = 7392.1113
$ws.Range("I84").Value = 8224.286
$ws.Range("K84").Value = 82242.86
$ws.Range("M84").Value = -76938.86
$ws.Range("H96").Value = 3832.5
$ws.Range("J96").Value = 4749
$ws.Range("L96").Value = 4749
$ws.Range("N96").Value = -7495
$ws.Range("H107").Value = 7240.75
$ws.Range("I107").Value = 11995.223
$ws.Range("J107").Value = 1127.8572
$ws.Range("K107").Value = 35985.669
$ws.Range("L107").Value = 3383.5716
$ws.Range("M107").Value = -34065.669
$ws.Range("N107").Value = -7223.571599999999
$ws.Range("H132").Value = 1472.4062
$ws.Range("I132").Value = 1307.6207
$ws.Range("K132").Value = 3922.8621
$ws.Range("M132").Value = -1392.8621
$ws.Range("H136").Value = 6343.24
$ws.Range("I136").Value = 5959.1777
$ws.Range("J136").Value = 9799.799999999999
$ws.Range("K136").Value = 17877.5331
$ws.Range("L136").Value = 29399.4
$ws.Range("M136").Value = -15327.5331
$ws.Range("N136").Value = -34499.39999999999
